$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old FAPs-prefixed rows (5:7) — their content becomes the
# new, recalculated rows 2:4.
$ws.Rows("5:7").Delete()

function Set-RowValues {
    param($ws, $row, $values)
    for ($i = 0; $i -lt $values.Length; $i++) {
        $ws.Cells.Item($row, $i + 1).Value = $values[$i]
    }
}

# Row 2: FAPs -> Wnt1 -> Fzd2 -> ECs
Set-RowValues $ws 2 @(
    "FAPs", "Wnt1", "Fzd2", "ECs",
    3, 1, 0.3532066666666667, 1.05962, 1, 1,
    2, 0.6666666666666666, 0.02648366666666667, 0.07945099999999999,
    0.001430039273477916, 0.001430039273477917, 0.009354207624444444, 0.08418786861999999,
    0.001430039273477916, 0.001430039273477917
)

# Row 3: FAPs -> Wnt1 -> Fzd2 -> FAPs
Set-RowValues $ws 3 @(
    "FAPs", "Wnt1", "Fzd2", "FAPs",
    3, 1, 0.3532066666666667, 1.05962, 1, 1,
    3, 1, 12.99468133333333, 38.984044,
    0.7016741634339546, 0.7016741634339547, 4.589808078142222, 41.30827270328,
    0.7016741634339546, 0.7016741634339547
)

# Row 4: FAPs -> Wnt1 -> Fzd2 -> MuSCs
Set-RowValues $ws 4 @(
    "FAPs", "Wnt1", "Fzd2", "MuSCs",
    3, 1, 0.3532066666666667, 1.05962, 1, 1,
    3, 1, 5.498373, 16.495119,
    0.2968957972925674, 0.2968957972925675, 1.94206199942, 17.47855799478,
    0.2968957972925674, 0.2968957972925675
)

Write-Host "Done"
